# Updates the "KEIPI_CCSSMR" sheet with new CCS efficiency / captured CO2 data.
#
# Summary of the edit:
#  1. Insert a new row above the "Output conversion factor" row to hold a
#     new "CCS Efficiency" input (0.9).
#  2. Add a "Captured CO2" line that multiplies the existing Emissions value
#     by the H2 output and the new CCS Efficiency factor.
#  3. Make the CCS-ELC energy penalty (row with "CCSELC") a live formula that
#     pulls from the recomputed "Captured CO2" chain, instead of a frozen
#     number.
#  4. Highlight the resulting GWh total in red, like the rest of that
#     computed total row.
#  5. Leave the active-cell selection on F18, matching where the author was
#     last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KEIPI_CCSSMR")

# 1. Insert new row 5 ("CCS Efficiency" input). Excel shifts every row/
#    formula reference at or below row 5 down by one automatically.
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "CCS Efficiency"
$ws.Range("B5").Value = 0.9

# 2. New "Captured CO2" row (this is the old "=B33*F8" emissions-total row,
#    now shifted to row 34, with a new row 35 appended directly below it).
$ws.Range("A35").Value = "Captured CO2"
$ws.Range("A35").Font.Bold = $true
$ws.Range("B35").Formula = "=B34*F9*B5"

# 3. CCS ELC penalty (row 17 after the shift) now references the bottom
#    calculation (B39) instead of being a hard-coded number.
$ws.Range("F17").Formula = "=B39"

# 4. Match the red font styling already used for this kind of output total.
$ws.Range("B39").Font.Color = 255
$ws.Range("C39").Font.Color = 255

# 5. Restore the author's last-used selection.
$null = $ws.Range("F18").Select()
